# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# Source diff updates Price (D) / Volume(1h) (E) text values for most rows,
# and swaps the Toncoin / WEMIXToken rows (26 <-> 27) including their links.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.039.26'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').Value = '2.276.89'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').Value = '''231.21'
$ws.Range('E5').Value = '  -0.49%  '

$ws.Range('E6').Value = '  +1.93%  '

$ws.Range('D7').Value = '''63.88'
$ws.Range('E7').Value = '  +4.01%  '

$ws.Range('D9').Value = '''0.448'
$ws.Range('E9').Value = '  +8.89%  '

$ws.Range('E10').Value = '  +10.19%  '

$ws.Range('D11').Value = '''57.04'
$ws.Range('E11').Value = '  -1.46%  '

$ws.Range('D12').Value = '''26.77'
$ws.Range('E12').Value = '  +18.27%  '

$ws.Range('E13').Value = '  +2.21%  '

$ws.Range('D14').Value = '2.617.78'
$ws.Range('E14').Value = '  +0.84%  '

$ws.Range('D15').Value = '''15.73'
$ws.Range('E15').Value = '  -0.05%  '

$ws.Range('D16').Value = '''6.15'
$ws.Range('E16').Value = '  +7.78%  '

$ws.Range('D17').Value = '''0.842'
$ws.Range('E17').Value = '  +3.70%  '

$ws.Range('D18').Value = '2.278.52'
$ws.Range('E18').Value = '  +1.21%  '

$ws.Range('D19').Value = '43.960.07'
$ws.Range('E19').Value = '  +1.43%  '

$ws.Range('E20').Value = '  +7.06%  '

$ws.Range('D21').Value = '''73.81'
$ws.Range('E21').Value = '  +1.02%  '

$ws.Range('D22').Value = '''6.12'
$ws.Range('E22').Value = '  -1.72%  '

$ws.Range('D23').Value = '''253.61'
$ws.Range('E23').Value = '  +2.17%  '

$ws.Range('E24').Value = '  +0.28%  '

$ws.Range('E25').Value = '  -5.25%  '

$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D26').Value = '''3.36'
$ws.Range('E26').Value = '  +25.81%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '''2.28'
$ws.Range('E27').Value = '  -5.16%  '

$ws.Range('D28').Value = '''10.06'
$ws.Range('E28').Value = '  +2.40%  '

$ws.Range('D29').Value = '''171.80'
$ws.Range('E29').Value = '  +1.29%  '

$ws.Range('E30').Value = '  -3.03%  '

$ws.Range('D31').Value = '''20.87'
$ws.Range('E31').Value = '  +1.68%  '

$ws.Range('E32').Value = '  -7.07%  '

$ws.Range('E33').Value = '  +3.30%  '

$ws.Range('D34').Value = '''0.0701'
$ws.Range('E34').Value = '  +6.49%  '

$ws.Range('D35').Value = '''4.83'
$ws.Range('E35').Value = '  +1.90%  '

$ws.Range('D36').Value = '''4.91'
$ws.Range('E36').Value = '  -2.79%  '

$ws.Range('D37').Value = '''3.79'
$ws.Range('E37').Value = '  +4.80%  '

$ws.Range('D38').Value = '''6.55'
$ws.Range('E38').Value = '  +0.83%  '

$ws.Range('E39').Value = '  -3.27%  '

$ws.Range('D40').Value = '''0.0260'
$ws.Range('E40').Value = '  +3.52%  '

$ws.Range('D41').Value = '''0.000241'
$ws.Range('E41').Value = '  +5.75%  '

$ws.Range('E42').Value = '  +0.32%  '

$ws.Range('D43').Value = '''17.68'
$ws.Range('E43').Value = '  +5.87%  '

$ws.Range('E44').Value = '  +0.01%  '

$ws.Range('D45').Value = '''8.28'
$ws.Range('E45').Value = '  -4.95%  '

$ws.Range('D46').Value = '''10.43'
$ws.Range('E46').Value = '  +15.22%  '

$ws.Range('D47').Value = '''98.59'
$ws.Range('E47').Value = '  +1.17%  '

$ws.Range('E48').Value = '  -1.25%  '

$ws.Range('D49').Value = '''4.36'
$ws.Range('E49').Value = '  -1.71%  '

$ws.Range('D50').Value = '1.451.21'
$ws.Range('E50').Value = '  -1.14%  '

$ws.Range('D51').Value = '''2.30'
$ws.Range('E51').Value = '  +2.68%  '
